$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 2.62
$ws.Range("I2").Value  = 2.9
$ws.Range("J2").Value  = 3.25
$ws.Range("L2").Value  = 3.55
$ws.Range("V2").Value  = 1.9
$ws.Range("W2").Value  = 7.5
$ws.Range("Y2").Value  = 9.5
$ws.Range("Z2").Value  = 32
$ws.Range("AA2").Value = 24
$ws.Range("AB2").Value = 32
$ws.Range("AF2").Value = 70
$ws.Range("AH2").Value = 7.7
$ws.Range("AI2").Value = 14.5
$ws.Range("AL2").Value = 28
$ws.Range("AM2").Value = 37
$ws.Range("AN2").Value = 4.55
$ws.Range("AO2").Value = 14.5
$ws.Range("AP2").Value = 21
$ws.Range("AQ2").Value = 65
$ws.Range("AR2").Value = 90
$ws.Range("AX2").Value = 17
$ws.Range("AY2").Value = 23
$ws.Range("BA2").Value = 120
